$d = $word.ActiveDocument

$replacements = @(
    @("2024-06-09 Sunday", "2024-06-10 Monday"),
    @("61×33=", "29×46="),
    @("29×82=", "30×18="),
    @("38×59=", "89×52="),
    @("18×59=", "81×69="),
    @("60×37=", "33×82="),
    @("77×94=", "58×49="),
    @("28×71=", "33×45="),
    @("82×91=", "14×28="),
    @("64×33=", "25×11="),
    @("33×76=", "44×27="),
    @("36×54=", "29×40="),
    @("31×28=", "92×99="),
    @("48×16=", "51×17="),
    @("30×63=", "13×34="),
    @("28×38=", "74×28="),
    @("95×72=", "27×58="),
    @("31×67=", "49×85="),
    @("68×95=", "41×48="),
    @("19×28=", "58×48="),
    @("41×49=", "56×74="),
    @("22×45=", "82×45="),
    @("16×19=", "37×50="),
    @("94×56=", "99×36="),
    @("65×61=", "58×74="),
    @("59×46=", "31×75=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
